$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename tour result files to reflect the ".tour" output extension
$ws.Range("A2").Value = "tsp_example_1.txt.tour"
$ws.Range("A3").Value = "tsp_example_2.txt.tour"
$ws.Range("A4").Value = "tsp_example_3.txt.tour"

# Updated timing measurement for the first example
$ws.Range("C2").Value = 0.059

# Give the report table a thin box border around every cell
$ws.Range("A1:C4").Borders.LineStyle = 1

# Widen column A so the longer filenames (with the .tour suffix) are readable
$ws.Columns.Item(1).ColumnWidth = 22.8333333333

# Leave the selection where the user last clicked while reviewing the report
$ws.Range("A11").Select() | Out-Null
